$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.563.13'
$ws.Range('E2').Value = '  -1.59%  '
$ws.Range('D3').Value = '2.433.43'
$ws.Range('E3').Value = '  -2.10%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '514.68'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.32%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '130.10'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.94%  '
$ws.Range('E7').Value = '  -0.16%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.550'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.98%  '
$ws.Range('D9').Value = '2.443.84'
$ws.Range('E9').Value = '  -1.78%  '
$ws.Range('E10').Value = '  -0.09%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0951'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -5.22%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.18'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.92%  '
$ws.Range('E13').Value = '  -3.67%  '
$ws.Range('D14').Value = '2.867.07'
$ws.Range('E14').Value = '  -2.10%  '
$ws.Range('D15').Value = '57.470.74'
$ws.Range('E15').Value = '  -1.60%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '21.78'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E17').Value = '  -3.33%  '
$ws.Range('D18').Value = '2.437.57'
$ws.Range('E18').Value = '  -2.05%  '
$ws.Range('E19').Value = '  -4.46%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '315.98'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.80%  '
$ws.Range('E21').Value = '  -2.67%  '
$ws.Range('E22').Value = '  +0.08%  '
$ws.Range('E23').Value = '  -2.72%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '63.42'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.47%  '
$ws.Range('E25').Value = '  -1.65%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.997'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.12%  '
$ws.Range('E28').Value = '  -3.47%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '170.45'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +3.23%  '
$ws.Range('E30').Value = '  -4.06%  '
$ws.Range('E31').Value = '  -2.62%  '
$ws.Range('E32').Value = '  -2.72%  '
$ws.Range('E33').Value = '  +2.39%  '
$ws.Range('E35').Value = '  -0.20%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '17.68'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.31%  '
$ws.Range('E37').Value = '  -4.70%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.92'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.88%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '36.22'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.34%  '
$ws.Range('E40').Value = '  -2.98%  '
$ws.Range('E41').Value = '  -2.65%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '271.12'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.78%  '
$ws.Range('B43').Value = 'Filecoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.38'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -4.49%  '
$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.92'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.96%  '
$ws.Range('E45').Value = '  -2.02%  '
$ws.Range('E46').Value = '  -0.92%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '120.58'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -5.40%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0484'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.37%  '
$ws.Range('E49').Value = '  -3.20%  '
$ws.Range('E50').Value = '  -3.84%  '
$ws.Range('D51').Value = '1.707.91'
$ws.Range('E51').Value = '  -1.91%  '
